# Update LR-pair data for Sertad1-Ar with new TPM-based values.
# Data now covers the full 4x4 sending/target cluster grid (rows 2-17).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @(
    @{rowNum=2; colA="ECs"; colD="ECs"; colE=3; colF=1; colG=12.57438866666667; colH=37.723166; colI=0.2319227286520073; colJ=0.2319227286520073; colK=3; colL=1; colM=0.9379940000000001; colN=2.813982; colO=0.0640701907773854; colP=0.0640701907773854; colQ=11.79470112300133; colR=106.152310107012; colS=0.0148593334703459; colT=0.0148593334703459}
    @{rowNum=3; colA="ECs"; colD="FAPs"; colE=3; colF=1; colG=12.57438866666667; colH=37.723166; colI=0.2319227286520073; colJ=0.2319227286520073; colK=3; colL=1; colM=8.022254333333334; colN=24.066763; colO=0.5479644492410116; colP=0.5479644492410115; colQ=100.8749439701842; colR=907.874495731658; colS=0.1270854102722698; colT=0.1270854102722697}
    @{rowNum=4; colA="ECs"; colD="MuSCs"; colE=3; colF=1; colG=12.57438866666667; colH=37.723166; colI=0.2319227286520073; colJ=0.2319227286520073; colK=3; colL=1; colM=5.642847666666667; colN=16.928543; colO=0.3854377816180673; colP=0.3854377816180672; colQ=70.95535974745978; colR=638.598237727138; colS=0.08939178203843867; colT=0.08939178203843866}
    @{rowNum=5; colA="ECs"; colD="Resolving-Mac"; colE=3; colF=1; colG=12.57438866666667; colH=37.723166; colI=0.2319227286520073; colJ=0.2319227286520073; colK=1; colL=0.3333333333333333; colM=0.037004; colN=0.111012; colO=0.002527578363535768; colP=0.002527578363535768; colQ=0.4653026782213333; colR=4.187724103992; colS=0.0005862028709529907; colT=0.0005862028709529906}
    @{rowNum=6; colA="FAPs"; colD="ECs"; colE=3; colF=1; colG=16.174732; colH=48.524196; colI=0.2983276626878248; colJ=0.2983276626878247; colK=3; colL=1; colM=0.9379940000000001; colN=2.813982; colO=0.0640701907773854; colP=0.0640701907773854; colQ=15.171801567608; colR=136.546214108472; colS=0.01911391026258042; colT=0.01911391026258041}
    @{rowNum=7; colA="FAPs"; colD="FAPs"; colE=3; colF=1; colG=16.174732; colH=48.524196; colI=0.2983276626878248; colJ=0.2983276626878247; colK=3; colL=1; colM=8.022254333333334; colN=24.066763; colO=0.5479644492410116; colP=0.5479644492410115; colQ=129.7578138775054; colR=1167.820324897548; colS=0.1634729533780922; colT=0.1634729533780921}
    @{rowNum=8; colA="FAPs"; colD="MuSCs"; colE=3; colF=1; colG=16.174732; colH=48.524196; colI=0.2983276626878248; colJ=0.2983276626878247; colK=3; colL=1; colM=5.642847666666667; colN=16.928543; colO=0.3854377816180673; colP=0.3854377816180672; colQ=91.2715487251587; colR=821.4439385264282; colS=0.1149867525016983; colT=0.1149867525016982}
    @{rowNum=9; colA="FAPs"; colD="Resolving-Mac"; colE=3; colF=1; colG=16.174732; colH=48.524196; colI=0.2983276626878248; colJ=0.2983276626878247; colK=1; colL=0.3333333333333333; colM=0.037004; colN=0.111012; colO=0.002527578363535768; colP=0.002527578363535768; colQ=0.5985297829280001; colR=5.386768046352; colS=0.0007540465454539428; colT=0.0007540465454539426}
    @{rowNum=10; colA="MuSCs"; colD="ECs"; colE=3; colF=1; colG=13.006622; colH=39.019866; colI=0.2398948644542636; colJ=0.2398948644542636; colK=3; colL=1; colM=0.9379940000000001; colN=2.813982; colO=0.0640701907773854; colP=0.0640701907773854; colQ=12.200133396268; colR=109.801200566412; colS=0.01537010973209968; colT=0.01537010973209968}
    @{rowNum=11; colA="MuSCs"; colD="FAPs"; colE=3; colF=1; colG=13.006622; colH=39.019866; colI=0.2398948644542636; colJ=0.2398948644542636; colK=3; colL=1; colM=8.022254333333334; colN=24.066763; colO=0.5479644492410116; colP=0.5479644492410115; colQ=104.3424297015287; colR=939.081867313758; colS=0.1314538572764277; colT=0.1314538572764277}
    @{rowNum=12; colA="MuSCs"; colD="MuSCs"; colE=3; colF=1; colG=13.006622; colH=39.019866; colI=0.2398948644542636; colJ=0.2398948644542636; colK=3; colL=1; colM=5.642847666666667; colN=16.928543; colO=0.3854377816180673; colP=0.3854377816180672; colQ=73.39438660391534; colR=660.549479435238; colS=0.09246454437681831; colT=0.0924645443768183}
    @{rowNum=13; colA="MuSCs"; colD="Resolving-Mac"; colE=3; colF=1; colG=13.006622; colH=39.019866; colI=0.2398948644542636; colJ=0.2398948644542636; colK=1; colL=0.3333333333333333; colM=0.037004; colN=0.111012; colO=0.002527578363535768; colP=0.002527578363535768; colQ=0.481297040488; colR=4.331673364392; colS=0.0006063530689179425; colT=0.0006063530689179424}
    @{rowNum=14; colA="Resolving-Mac"; colD="ECs"; colE=3; colF=1; colG=12.46226666666667; colH=37.3868; colI=0.2298547442059043; colJ=0.2298547442059043; colK=3; colL=1; colM=0.9379940000000001; colN=2.813982; colO=0.0640701907773854; colP=0.0640701907773854; colQ=11.68953135973334; colR=105.2057822376; colS=0.01472683731235941; colT=0.01472683731235941}
    @{rowNum=15; colA="Resolving-Mac"; colD="FAPs"; colE=3; colF=1; colG=12.46226666666667; colH=37.3868; colI=0.2298547442059043; colJ=0.2298547442059043; colK=3; colL=1; colM=8.022254333333334; colN=24.066763; colO=0.5479644492410116; colP=0.5479644492410115; colQ=99.97547276982223; colR=899.7792549284001; colS=0.125952228314222; colT=0.125952228314222}
    @{rowNum=16; colA="Resolving-Mac"; colD="MuSCs"; colE=3; colF=1; colG=12.46226666666667; colH=37.3868; colI=0.2298547442059043; colJ=0.2298547442059043; colK=3; colL=1; colM=5.642847666666667; colN=16.928543; colO=0.3854377816180673; colP=0.3854377816180672; colQ=70.32267238137779; colR=632.9040514324; colS=0.08859470270111207; colT=0.08859470270111205}
    @{rowNum=17; colA="Resolving-Mac"; colD="Resolving-Mac"; colE=3; colF=1; colG=12.46226666666667; colH=37.3868; colI=0.2298547442059043; colJ=0.2298547442059043; colK=1; colL=0.3333333333333333; colM=0.037004; colN=0.111012; colO=0.002527578363535768; colP=0.002527578363535768; colQ=0.4611537157333334; colR=4.1503834416; colS=0.0005809758782108922; colT=0.0005809758782108922}
)

foreach ($row in $rowData) {
    $rn = $row.rowNum
    $ws.Range("A$rn").Value = $row.colA
    $ws.Range("B$rn").Value = "Sertad1"
    $ws.Range("C$rn").Value = "Ar"
    $ws.Range("D$rn").Value = $row.colD
    $ws.Range("E$rn").Value = $row.colE
    $ws.Range("F$rn").Value = $row.colF
    $ws.Range("G$rn").Value = $row.colG
    $ws.Range("H$rn").Value = $row.colH
    $ws.Range("I$rn").Value = $row.colI
    $ws.Range("J$rn").Value = $row.colJ
    $ws.Range("K$rn").Value = $row.colK
    $ws.Range("L$rn").Value = $row.colL
    $ws.Range("M$rn").Value = $row.colM
    $ws.Range("N$rn").Value = $row.colN
    $ws.Range("O$rn").Value = $row.colO
    $ws.Range("P$rn").Value = $row.colP
    $ws.Range("Q$rn").Value = $row.colQ
    $ws.Range("R$rn").Value = $row.colR
    $ws.Range("S$rn").Value = $row.colS
    $ws.Range("T$rn").Value = $row.colT
}
